# Actualizacion TASK, Schedule, Week Personal
# Fill in the "Actual" Hours (K) and Week No. (M) columns for the task rows
# 32-38 on the TASK sheet. The "Cumulative Hours" column (L) already holds a
# running-sum formula (L[n] = K[n] + L[n-1]) so it recalculates automatically
# once K is populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 32
$ws.Range("K32").Value = 0.5
$ws.Range("M32").Value = 9

# Row 33
$ws.Range("K33").Value = 1
$ws.Range("M33").Value = 9

# Row 34
$ws.Range("K34").Value = 2
$ws.Range("M34").Value = 9

# Row 35
$ws.Range("K35").Value = 2.5
$ws.Range("M35").Value = 9

# Row 36
$ws.Range("K36").Value = 1.5
$ws.Range("M36").Value = 9

# Row 37
$ws.Range("K37").Value = 1.5
$ws.Range("M37").Value = 9

# Row 38
$ws.Range("K38").Value = 1.5
$ws.Range("M38").Value = 9

# Move the active selection / cursor position to reflect where the author
# left off editing (K38), matching the updated sheetView selection.
$ws.Range("K38").Select()
